# "Creado grafico de tipos de modelo"
# Add a new "MAE" metric column between "R2" and "Tipo", and refresh the
# MSE/R2 values plus the new MAE value for the existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Tipo" column (D) one place to the right, freeing up
# column D for the new "MAE" metric.
$ws.Columns.Item(4).Insert()

# New column header + value.
$ws.Range("D1").Value = "MAE"
$ws.Range("D2").Value = 0.1661107093174136

# Updated metric values for the existing row.
$ws.Range("B2").Value = 0.05090842587528387
$ws.Range("C2").Value = 0.9985028700472656
